$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NATMI LR-pairs data rows (2-7) with the newly re-run TPM values.
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric statistics recomputed from the refreshed TPM input.

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Ntf5"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.563831
$ws.Range("H2").Value = 1.691493
$ws.Range("I2").Value = 0.5426742997499485
$ws.Range("J2").Value = 0.5426742997499486
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.220618
$ws.Range("N2").Value = 0.6618539999999999
$ws.Range("O2").Value = 0.6621850925462731
$ws.Range("P2").Value = 0.6621850925462731
$ws.Range("Q2").Value = 0.124391267558
$ws.Range("R2").Value = 1.119521408022
$ws.Range("S2").Value = 0.3593508314024036
$ws.Range("T2").Value = 0.3593508314024037
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntf5"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.563831
$ws.Range("H3").Value = 1.691493
$ws.Range("I3").Value = 0.5426742997499485
$ws.Range("J3").Value = 0.5426742997499486
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1125486666666667
$ws.Range("N3").Value = 0.337646
$ws.Range("O3").Value = 0.3378149074537269
$ws.Range("P3").Value = 0.3378149074537269
$ws.Range("Q3").Value = 0.06345842727533334
$ws.Range("R3").Value = 0.5711258454779999
$ws.Range("S3").Value = 0.1833234683475449
$ws.Range("T3").Value = 0.1833234683475449
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ntf5"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05363299999999999
$ws.Range("H4").Value = 0.160899
$ws.Range("I4").Value = 0.05162052231695133
$ws.Range("J4").Value = 0.05162052231695134
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.220618
$ws.Range("N4").Value = 0.6618539999999999
$ws.Range("O4").Value = 0.6621850925462731
$ws.Range("P4").Value = 0.6621850925462731
$ws.Range("Q4").Value = 0.011832405194
$ws.Range("R4").Value = 0.106491646746
$ws.Range("S4").Value = 0.03418234034773737
$ws.Range("T4").Value = 0.03418234034773738
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Ntf5"
$ws.Range("C5").Value = "Ntrk1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05363299999999999
$ws.Range("H5").Value = 0.160899
$ws.Range("I5").Value = 0.05162052231695133
$ws.Range("J5").Value = 0.05162052231695134
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1125486666666667
$ws.Range("N5").Value = 0.337646
$ws.Range("O5").Value = 0.3378149074537269
$ws.Range("P5").Value = 0.3378149074537269
$ws.Range("Q5").Value = 0.006036322639333333
$ws.Range("R5").Value = 0.054326903754
$ws.Range("S5").Value = 0.01743818196921396
$ws.Range("T5").Value = 0.01743818196921396
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Ntf5"
$ws.Range("C6").Value = "Ntrk1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.421522
$ws.Range("H6").Value = 1.264566
$ws.Range("I6").Value = 0.4057051779331001
$ws.Range("J6").Value = 0.4057051779331002
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.220618
$ws.Range("N6").Value = 0.6618539999999999
$ws.Range("O6").Value = 0.6621850925462731
$ws.Range("P6").Value = 0.6621850925462731
$ws.Range("Q6").Value = 0.09299534059599998
$ws.Range("R6").Value = 0.8369580653639999
$ws.Range("S6").Value = 0.2686519207961321
$ws.Range("T6").Value = 0.2686519207961321
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Ntf5"
$ws.Range("C7").Value = "Ntrk1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.421522
$ws.Range("H7").Value = 1.264566
$ws.Range("I7").Value = 0.4057051779331001
$ws.Range("J7").Value = 0.4057051779331002
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1125486666666667
$ws.Range("N7").Value = 0.337646
$ws.Range("O7").Value = 0.3378149074537269
$ws.Range("P7").Value = 0.3378149074537269
$ws.Range("Q7").Value = 0.04744173907066666
$ws.Range("R7").Value = 0.4269756516359999
$ws.Range("S7").Value = 0.137053257136968
$ws.Range("T7").Value = 0.137053257136968
